$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cost burden values in row 2
$ws.Range("O2").Value = 847.32853333333333
$ws.Range("P2").Value = 141.43866666666668
$ws.Range("Q2").Value = 784.93133333333333
$ws.Range("R2").Value = 159.58693333333335
$ws.Range("S2").Value = 10062.020867460158

# Update the selection to A2:T2 with active cell A2
$ws.Range("A2:T2").Select()
